# Weekly price-sheet update for "Cebollín baby" (Agrícola del Norte S.A. de Arica).
# A new week's record is inserted as row 122 (pushing the existing rows 122-144
# down to 123-145); only the price/date columns differ between records, the
# market/category descriptive columns repeat unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 122, shifting 122-144 -> 123-145.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new weekly record.
$ws.Cells.Item(122, 1).Value  = 1
$ws.Cells.Item(122, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(122, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(122, 4).Value  = 45209
$ws.Cells.Item(122, 5).Value  = 15
$ws.Cells.Item(122, 6).Value  = 100112038
$ws.Cells.Item(122, 7).Value  = "Cebollín baby"
$ws.Cells.Item(122, 8).Value  = "Sin especificar"
$ws.Cells.Item(122, 9).Value  = "Primera"
$ws.Cells.Item(122, 10).Value = 220
$ws.Cells.Item(122, 11).Value = 2000
$ws.Cells.Item(122, 12).Value = 2000
$ws.Cells.Item(122, 13).Value = 2000
$ws.Cells.Item(122, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 1000
$ws.Cells.Item(122, 17).Value = 2
$ws.Cells.Item(122, 18).Value = "Hortaliza"
